$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H1").Value = "reference_period"
$ws.Range("I1").Value = "remarks"
$ws.Range("H2").Value = 2020
$ws.Range("I2").Value = "Test note"
Write-Host "done"
